# Applies the "fac-template" edit described by the commit diff:
#  - H19:H44 get "=SUM(G{row}*F{row})" line-total formulas
#  - B46 gets "=SUM(F19:F44)" (total quantity) and is centered
#  - E46 gets "=SUM(H47-H45)" (same value as H46), formatted like the
#    existing currency column and right aligned
#  - F19:F44 (the "Qté" column) becomes horizontally centered

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Line-total formulas for each article row (H = G * F)
for ($row = 19; $row -le 44; $row++) {
    $ws.Range("H$row").Formula = "=SUM(G$row*F$row)"
}

# 2) Center the quantity column ("Qté", F19:F44)
$ws.Range("F19:F44").HorizontalAlignment = -4108

# 3) Total quantity formula in B46, centered like the header style above it
$ws.Range("B46").Formula = "=SUM(F19:F44)"
$ws.Range("B46").HorizontalAlignment = -4108

# 4) E46 mirrors H46 ("TVA amount") - same formula, currency number format,
#    copied border/font from its left neighbour D46, right aligned
$ws.Range("D46").Copy()
$ws.Range("E46").PasteSpecial(-4122)
$ws.Range("E46").Formula = "=SUM(H47-H45)"
$ws.Range("E46").NumberFormat = $ws.Range("H45").NumberFormat
$ws.Range("E46").HorizontalAlignment = -4152
